{"js": "// Apply text replacements describing the diff: update the date and the 25\n// division-problem table cells. Each (old, new) pair is unique in the document,\n// so a body.search() + insertText(replace) per pair is safe and preserves run formatting.\nconst replacements = [\n  [\"2026-02-27 Friday\", \"2026-02-28 Saturday\"],\n  [\"412\u00f73=137, 1\", \"855\u00f72=427, 1\"],\n  [\"749\u00f78=93, 5\", \"739\u00f79=82, 1\"],\n  [\"385\u00f74=96, 1\", \"140\u00f75=28, 0\"],\n  [\"244\u00f75=48, 4\", \"681\u00f79=75, 6\"],\n  [\"708\u00f77=101, 1\", \"896\u00f76=149, 2\"],\n  [\"564\u00f78=70, 4\", \"101\u00f76=16, 5\"],\n  [\"980\u00f73=326, 2\", \"213\u00f79=23, 6\"],\n  [\"180\u00f76=30, 0\", \"107\u00f77=15, 2\"],\n  [\"689\u00f73=229, 2\", \"176\u00f76=29, 2\"],\n  [\"192\u00f74=48, 0\", \"651\u00f79=72, 3\"],\n  [\"620\u00f72=310, 0\", \"660\u00f73=220, 0\"],\n  [\"793\u00f78=99, 1\", \"349\u00f72=174, 1\"],\n  [\"649\u00f73=216, 1\", \"332\u00f75=66, 2\"],\n  [\"903\u00f73=301, 0\", \"618\u00f78=77, 2\"],\n  [\"848\u00f73=282, 2\", \"104\u00f72=52, 0\"],\n  [\"454\u00f76=75, 4\", \"450\u00f72=225, 0\"],\n  [\"865\u00f78=108, 1\", \"872\u00f74=218, 0\"],\n  [\"384\u00f77=54, 6\", \"189\u00f76=31, 3\"],\n  [\"202\u00f78=25, 2\", \"104\u00f79=11, 5\"],\n  [\"976\u00f75=195, 1\", \"713\u00f73=237, 2\"],\n  [\"101\u00f75=20, 1\", \"679\u00f78=84, 7\"],\n  [\"668\u00f72=334, 0\", \"277\u00f74=69, 1\"],\n  [\"402\u00f78=50, 2\", \"139\u00f73=46, 1\"],\n  [\"467\u00f73=155, 2\", \"493\u00f74=123, 1\"],\n  [\"516\u00f74=129, 0\", \"555\u00f74=138, 3\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# (old, new) text pairs describing the diff: the date line plus the 25 unique\n# division-problem cells. Each old string occurs exactly once in the document,\n# so Find/Replace per pair (wrapped with whole-story range + wildcard off) is safe.\n$pairs = @(\n    ,@(\"2026-02-27 Friday\", \"2026-02-28 Saturday\")\n    ,@(\"412\u00f73=137, 1\", \"855\u00f72=427, 1\")\n    ,@(\"749\u00f78=93, 5\", \"739\u00f79=82, 1\")\n    ,@(\"385\u00f74=96, 1\", \"140\u00f75=28, 0\")\n    ,@(\"244\u00f75=48, 4\", \"681\u00f79=75, 6\")\n    ,@(\"708\u00f77=101, 1\", \"896\u00f76=149, 2\")\n    ,@(\"564\u00f78=70, 4\", \"101\u00f76=16, 5\")\n    ,@(\"980\u00f73=326, 2\", \"213\u00f79=23, 6\")\n    ,@(\"180\u00f76=30, 0\", \"107\u00f77=15, 2\")\n    ,@(\"689\u00f73=229, 2\", \"176\u00f76=29, 2\")\n    ,@(\"192\u00f74=48, 0\", \"651\u00f79=72, 3\")\n    ,@(\"620\u00f72=310, 0\", \"660\u00f73=220, 0\")\n    ,@(\"793\u00f78=99, 1\", \"349\u00f72=174, 1\")\n    ,@(\"649\u00f73=216, 1\", \"332\u00f75=66, 2\")\n    ,@(\"903\u00f73=301, 0\", \"618\u00f78=77, 2\")\n    ,@(\"848\u00f73=282, 2\", \"104\u00f72=52, 0\")\n    ,@(\"454\u00f76=75, 4\", \"450\u00f72=225, 0\")\n    ,@(\"865\u00f78=108, 1\", \"872\u00f74=218, 0\")\n    ,@(\"384\u00f77=54, 6\", \"189\u00f76=31, 3\")\n    ,@(\"202\u00f78=25, 2\", \"104\u00f79=11, 5\")\n    ,@(\"976\u00f75=195, 1\", \"713\u00f73=237, 2\")\n    ,@(\"101\u00f75=20, 1\", \"679\u00f78=84, 7\")\n    ,@(\"668\u00f72=334, 0\", \"277\u00f74=69, 1\")\n    ,@(\"402\u00f78=50, 2\", \"139\u00f73=46, 1\")\n    ,@(\"467\u00f73=155, 2\", \"493\u00f74=123, 1\")\n    ,@(\"516\u00f74=129, 0\", \"555\u00f74=138, 3\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"Text not found: $old\"\n    }\n}"}
